# Update the "Price" column (D) with refreshed values, matching the
# "Updated symbol list ... with GitHub Actions" commit. Values are stored
# as text (they include placeholders like "--" elsewhere in the column),
# so each new value is entered with a leading apostrophe to force Excel
# to keep it as text instead of auto-converting to a number (which would
# also silently drop meaningful trailing zeros, e.g. "5.390" -> 5.39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "242.96"
    "D3"  = "23.17"
    "D4"  = "5.390"
    "D6"  = "3.397"
    "D7"  = "0.8065"
    "D8"  = "0.9095"
    "D10" = "0.07437"
    "D11" = "0.03332"
    "D12" = "0.03041"
    "D13" = "0.09329"
    "D14" = "3.948"
    "D15" = "0.001577"
    "D16" = "0.04782"
    "D17" = "0.0005945"
    "D18" = "0.006141"
    "D20" = "0.004415"
    "D21" = "0.0009876"
    "D22" = "0.00007806"
    "D40" = "0.03873"
    "D41" = "0.006200"
    "D42" = "0.1066"
    "D44" = "0.007246"
    "D45" = "0.00005185"
    "D47" = "0.0005805"
    "D48" = "0.9107"
    "D50" = "0.00002102"
    "D51" = "0.0002002"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
